$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "paiewise"

# "Sep-22-2023" must land as literal text, not get auto-converted to a date
# serial by Excel's smart input parsing. Enter it as a formula returning the
# string, then paste-special as values so the cell ends up a plain text
# (shared-string) cell with no special number formatting applied.
$ws.Range("B12").Formula = '="Sep-22-2023"'
$ws.Range("B12").Copy()
$ws.Range("B12").PasteSpecial(-4163)

$ws.Range("C12").Value = "NV"
$ws.Range("D12").Value = "all_submitted_tracker_ninaSep-22-2023"

$ws.Range("D14").Select()
